$d = $word.ActiveDocument

# Update the date heading (first paragraph of the document)
$d.Paragraphs(1).Range.Text = "2023-08-10 Thursday"

# Update the division problems in the table, cell by cell so each
# assignment only touches the single target cell (avoids accidental
# matches against identical text elsewhere in the table).
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "29÷8=3, 5"  # was: 67÷9=7, 4
$tbl.Cell(1, 2).Range.Text = "18÷3=6, 0"  # was: 62÷7=8, 6
$tbl.Cell(1, 3).Range.Text = "33÷3=11, 0"  # was: 82÷4=20, 2
$tbl.Cell(1, 4).Range.Text = "40÷6=6, 4"  # was: 82÷3=27, 1
$tbl.Cell(1, 5).Range.Text = "23÷8=2, 7"  # was: 10÷6=1, 4
$tbl.Cell(5, 1).Range.Text = "45÷7=6, 3"  # was: 12÷3=4, 0
$tbl.Cell(5, 2).Range.Text = "14÷3=4, 2"  # was: 29÷8=3, 5
$tbl.Cell(5, 3).Range.Text = "14÷3=4, 2"  # was: 26÷9=2, 8
$tbl.Cell(5, 4).Range.Text = "55÷5=11, 0"  # was: 67÷8=8, 3
$tbl.Cell(5, 5).Range.Text = "31÷6=5, 1"  # was: 69÷5=13, 4
$tbl.Cell(9, 1).Range.Text = "83÷9=9, 2"  # was: 14÷8=1, 6
$tbl.Cell(9, 2).Range.Text = "55÷4=13, 3"  # was: 58÷9=6, 4
$tbl.Cell(9, 3).Range.Text = "88÷2=44, 0"  # was: 60÷3=20, 0
$tbl.Cell(9, 4).Range.Text = "48÷2=24, 0"  # was: 85÷4=21, 1
$tbl.Cell(9, 5).Range.Text = "58÷2=29, 0"  # was: 95÷7=13, 4
$tbl.Cell(13, 1).Range.Text = "78÷5=15, 3"  # was: 90÷3=30, 0
$tbl.Cell(13, 2).Range.Text = "95÷3=31, 2"  # was: 25÷6=4, 1
$tbl.Cell(13, 3).Range.Text = "57÷3=19, 0"  # was: 31÷9=3, 4
$tbl.Cell(13, 4).Range.Text = "66÷3=22, 0"  # was: 35÷9=3, 8
$tbl.Cell(13, 5).Range.Text = "42÷7=6, 0"  # was: 46÷9=5, 1
$tbl.Cell(17, 1).Range.Text = "96÷2=48, 0"  # was: 64÷9=7, 1
$tbl.Cell(17, 2).Range.Text = "58÷2=29, 0"  # was: 24÷4=6, 0
$tbl.Cell(17, 3).Range.Text = "87÷4=21, 3"  # was: 26÷2=13, 0
$tbl.Cell(17, 4).Range.Text = "78÷5=15, 3"  # was: 10÷3=3, 1
$tbl.Cell(17, 5).Range.Text = "81÷2=40, 1"  # was: 97÷7=13, 6
